$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated ligand/receptor expression values (Natmi following Dr Hou advice)
# Columns: E,F,G,H,I,J,K,L,M,N,O,P,Q,R,S,T  (F and L unchanged = 1)
$data = @{
    2  = @(3, 23.40397433333333, 70.211923,          0.5512176233382776, 0.5512176233382777, 3, 18.444833,          55.33449900000001, 0.529296397589589,  0.5292963975895891, 431.6823981146197, 3885.141583031577, 0.2917575023208453,  0.2917575023208454)
    3  = @(3, 23.40397433333333, 70.211923,          0.5512176233382776, 0.5512176233382777, 3, 8.028767999999999, 24.086304,          0.2303950368909585, 0.2303950368909585, 187.905080200288,  1691.145721802592, 0.1269978046639689,  0.126997804663969)
    4  = @(3, 23.40397433333333, 70.211923,          0.5512176233382776, 0.5512176233382777, 3, 8.374233,          25.122699,          0.2403085655194523, 0.2403085655194524, 195.990334193353,  1763.913007740177, 0.1324623163534633,  0.1324623163534633)
    5  = @(3, 10.21452833333333, 30.643585,          0.240575722363629,  0.240575722363629,  3, 18.444833,          55.33449900000001, 0.529296397589589,  0.5292963975895891, 188.4052692821017, 1695.647423538915, 0.1273358631945819,  0.127335863194582)
    6  = @(3, 10.21452833333333, 30.643585,          0.240575722363629,  0.240575722363629,  3, 8.028767999999999, 24.086304,          0.2303950368909585, 0.2303950368909585, 82.01007821776,    738.09070395984,   0.05542745242903728, 0.0554274524290373)
    7  = @(3, 10.21452833333333, 30.643585,          0.240575722363629,  0.240575722363629,  3, 8.374233,          25.122699,          0.2403085655194523, 0.2403085655194524, 85.53884024843501, 769.8495622359151, 0.0578124067400097,  0.05781240674000972)
    8  = @(3, 8.840180333333334, 26.520541,          0.2082066542980933, 0.2082066542980934, 3, 18.444833,          55.33449900000001, 0.529296397589589,  0.5292963975895891, 163.0556499382177, 1467.500849443959, 0.1102030320741617,  0.1102030320741618)
    9  = @(3, 8.840180333333334, 26.520541,          0.2082066542980933, 0.2082066542980934, 3, 8.028767999999999, 24.086304,          0.2303950368909585, 0.2303950368909585, 70.97575697449601, 638.7818127704639, 0.04796977979795226, 0.04796977979795227)
    10 = @(3, 8.840180333333334, 26.520541,          0.2082066542980933, 0.2082066542980934, 3, 8.374233,          25.122699,          0.2403085655194523, 0.2403085655194524, 74.02972987335102, 666.267568860159,  0.05003384242597933, 0.05003384242597934)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 5).Value  = $vals[0]   # E - Ligand-expressing cells
    $ws.Cells.Item($row, 7).Value  = $vals[1]   # G - Ligand average expression value
    $ws.Cells.Item($row, 8).Value  = $vals[2]   # H - Ligand total expression value
    $ws.Cells.Item($row, 9).Value  = $vals[3]   # I - Ligand derived specificity (avg)
    $ws.Cells.Item($row, 10).Value = $vals[4]   # J - Ligand derived specificity (total)
    $ws.Cells.Item($row, 11).Value = $vals[5]   # K - Receptor-expressing cells
    $ws.Cells.Item($row, 13).Value = $vals[6]   # M - Receptor average expression value
    $ws.Cells.Item($row, 14).Value = $vals[7]   # N - Receptor total expression value
    $ws.Cells.Item($row, 15).Value = $vals[8]   # O - Receptor derived specificity (avg)
    $ws.Cells.Item($row, 16).Value = $vals[9]   # P - Receptor derived specificity (total)
    $ws.Cells.Item($row, 17).Value = $vals[10]  # Q - Edge average expression weight
    $ws.Cells.Item($row, 18).Value = $vals[11]  # R - Edge total expression weight
    $ws.Cells.Item($row, 19).Value = $vals[12]  # S - Edge average expression derived specificity
    $ws.Cells.Item($row, 20).Value = $vals[13]  # T - Edge total expression derived specificity
}
